# "Generate Report for Handoff"
# The localization status report moved from "In Translation" to
# "Ready for handoff": update the status text, bump the HO xliff / handoff
# timestamps, and refresh the column widths that auto-size to the new
# (longer) status text on every sheet that shows it.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -----------------
$overview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$overview.Range("F2").Value = "Ready for handoff"   # de-de status
$zhcn.Range("C2").Value     = "Ready for handoff"   # Status column
$dede.Range("C2").Value     = "Ready for handoff"   # Status column

# --- Timestamps refreshed by the new handoff generation --------------------
$overview.Range("G2").Value = "2016-09-01 20:45:42" # Latest HO Xliff Generate Date
$dede.Range("H2").Value     = "2016-09-01 20:45:42" # Latest Handoff Datetime (de-de)
$zhcn.Range("H2").Value     = "2016-09-01 20:45:38" # Latest Handoff Datetime (zh-cn)

# --- Column widths: the status columns grew to fit "Ready for handoff" -----
$newStatusWidth = 16.333333333333332
$overview.Columns("E:F").ColumnWidth = $newStatusWidth
$zhcn.Columns("C:C").ColumnWidth = $newStatusWidth
$dede.Columns("C:C").ColumnWidth = $newStatusWidth
